$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B2: "5$August$92" -> "$#$%$$%^^"
$ws.Range("B2").Value = "$#$%$$%^^"

# Update the selection shown in the sheet view
$ws.Range("B9").Select()
